# Delete row 310 ("「停止距離の計算」" post) from Sheet1.
# All subsequent rows (311-352) shift up by one (to 310-351),
# and the sheet dimension shrinks from A1:C352 to A1:C351.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(310).Delete()
